# "updated burndown chart and scrum board"
# Move two scrum-board cards forward a column and update the sheet's
# view/selection state to match where the author left off working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")
$ws.Activate()

# --- Scrum board card moves -------------------------------------------------
# "Fazer os Use Case Diagrams (Todos tem de fazer)" moves from the
# "To do" column (D5) into the "Done" column (G6).
$card1 = $ws.Range("D5").Text
$ws.Range("G6").Value = $card1
$ws.Range("D5").ClearContents()

# "Identificar 3 Design Patterns (Todos tem de fazer)" moves from the
# "To do" column (D6) into the "Doing" column (E7).
$card2 = $ws.Range("D6").Text
$ws.Range("E7").Value = $card2
$ws.Range("D6").ClearContents()

# --- View state --------------------------------------------------------------
# Zoom out to 55% and scroll so column C is leftmost, then leave the
# selection on E14 (below the board), matching where editing stopped.
$excel.ActiveWindow.Zoom = 55
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1

$ws.Range("E14").Select()
